$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 148
$ws.Range("A148").Value = 'Rahmonova Oysuluv Mehmonaliyevna'
$ws.Range("B148").Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C148").Value = 'AB5022316'
$ws.Range("D148").Value = '''744'
$ws.Range("E148").Value = 'Andijon viloyati'
$ws.Range("F148").Value = 'Baliqchi tumani'
$ws.Range("G148").Value = '''998956767474'
$ws.Range("H148").Value = '14-11-2024'
$ws.Range("I148").Value = '''+998957677474'

# Row 149
$ws.Range("A149").Value = 'Sobirova Ismigul Bahromjon qizi'
$ws.Range("B149").Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C149").Value = 'AD2393298'
$ws.Range("D149").Value = '''745'
$ws.Range("E149").Value = 'Fargona viloyati'
$ws.Range("F149").Value = 'Buvayda tumani'
$ws.Range("G149").Value = '''998912032308'
$ws.Range("H149").Value = '14-11-2024'
$ws.Range("I149").Value = '''+998912032308'

# Row 150
$ws.Range("A150").Value = 'Muminova Mavludaxon Mahmudovna'
$ws.Range("B150").Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C150").Value = 'AC2511867'
$ws.Range("D150").Value = '''746'
$ws.Range("E150").Value = 'Andijon viloyati'
$ws.Range("F150").Value = 'Andijon tuman'
$ws.Range("G150").Value = '''998905407059'
$ws.Range("H150").Value = '15-11-2024'
$ws.Range("I150").Value = '''+998905407059'

# Row 151
$ws.Range("A151").Value = 'Ahmedova Nilufar Mirzaazizovna'
$ws.Range("B151").Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C151").Value = 'AA7377260'
$ws.Range("D151").Value = '''747'
$ws.Range("E151").Value = 'Andijon viloyati'
$ws.Range("F151").Value = 'Shahrixon tuman'
$ws.Range("G151").Value = '''998500721276'
$ws.Range("H151").Value = '18-11-2024'

# Row 152
$ws.Range("A152").Value = 'Samiyeva Farida Xudoyberdi qizi'
$ws.Range("B152").Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C152").Value = 'AD7058036'
$ws.Range("D152").Value = '''748'
$ws.Range("E152").Value = 'Jizzax viloyati'
$ws.Range("F152").Value = 'Arnasoy tumani'
$ws.Range("G152").Value = '''998936072294'
$ws.Range("H152").Value = '18-11-2024'
$ws.Range("I152").Value = '''+998936072294'

# Row 153
$ws.Range("A153").Value = 'Sobirova Zarnigor Sobirovna'
$ws.Range("B153").Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 864 soatlik'
$ws.Range("C153").Value = 'AD3746455'
$ws.Range("D153").Value = '''749'
$ws.Range("E153").Value = 'Qashqadaryo viloyati'
$ws.Range("F153").Value = 'Qarshi tumani'
$ws.Range("G153").Value = '''998500072344'
$ws.Range("H153").Value = '18-11-2024'
$ws.Range("I153").Value = '''+998905180033'

# Row 154
$ws.Range("A154").Value = 'Otaqoziyeva Gulhayo Mahmudjon qizi'
$ws.Range("B154").Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 864 soatlik'
$ws.Range("C154").Value = 'AD1672467'
$ws.Range("D154").Value = '''750'
$ws.Range("E154").Value = 'Fargona viloyati'
$ws.Range("F154").Value = 'Buvayda tumani'
$ws.Range("G154").Value = '''998910596680'
$ws.Range("H154").Value = '19-11-2024'
$ws.Range("I154").Value = '''+998910596680'

# Row 155
$ws.Range("A155").Value = 'Rasilova Sevinch Ilhomboyevna'
$ws.Range("B155").Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C155").Value = 'AD6148106'
$ws.Range("D155").Value = '''751'
$ws.Range("E155").Value = 'Toshkent shahri'
$ws.Range("F155").Value = 'Yangihayot tumani'
$ws.Range("G155").Value = '''998770723424'
$ws.Range("H155").Value = '19-11-2024'
$ws.Range("I155").Value = '''+998770273424'

# Row 156
$ws.Range("A156").Value = 'Norqulova Muqaddas Abdumannon qizi'
$ws.Range("B156").Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C156").Value = 'AD9376541'
$ws.Range("D156").Value = '''752'
$ws.Range("E156").Value = 'Andijon viloyati'
$ws.Range("F156").Value = 'Shahrixon tuman'
$ws.Range("G156").Value = '''998902201013'
$ws.Range("H156").Value = '19-11-2024'
$ws.Range("I156").Value = '''+998902201013'

# Row 157
$ws.Range("A157").Value = 'TO''XSANOVA FERUZA SOBIROVNA'
$ws.Range("B157").Value = 'Defektologiya (logopediya) 576 soatlik'
$ws.Range("C157").Value = 'AD7550400'
$ws.Range("D157").Value = '''753'
$ws.Range("E157").Value = 'Navoiy viloyati'
$ws.Range("F157").Value = 'Qiziltepa tumani'
$ws.Range("G157").Value = '''998912508529'
$ws.Range("H157").Value = '19-11-2024'
$ws.Range("I157").Value = '''+998912508529'

# Row 158
$ws.Range("A158").Value = 'TO''XSANOVA FERUZA SOBIROVNA'
$ws.Range("B158").Value = 'Defektologiya (logopediya) 576 soatlik'
$ws.Range("C158").Value = 'AD7550400'
$ws.Range("D158").Value = '''753'
$ws.Range("E158").Value = 'Navoiy viloyati'
$ws.Range("F158").Value = 'Qiziltepa tumani'
$ws.Range("G158").Value = '''998912508529'
$ws.Range("H158").Value = '19-11-2024'
$ws.Range("I158").Value = '''+998912508529'

# Row 159
$ws.Range("A159").Value = 'Quvondiqova Huriyat Bahronovna'
$ws.Range("B159").Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C159").Value = 'AD2409734'
$ws.Range("D159").Value = '''753'
$ws.Range("E159").Value = 'Toshkent shahri'
$ws.Range("F159").Value = 'Mirobod tumani'
$ws.Range("G159").Value = '''+998950038686'
$ws.Range("H159").Value = '19-11-2024'
$ws.Range("I159").Value = '''+998950038686'

# Row 160
$ws.Range("A160").Value = 'Hasanova Sevara ABDURAIMOVNA'
$ws.Range("B160").Value = 'Defektologiya (logopediya) 576 soatlik'
$ws.Range("C160").Value = 'AA9027122'
$ws.Range("D160").Value = '''753'
$ws.Range("E160").Value = 'Surxondaryo viloyati'
$ws.Range("F160").Value = 'Uzun tumani'
$ws.Range("G160").Value = '''998916158686'
$ws.Range("H160").Value = '19-11-2024'
$ws.Range("I160").Value = '''+998904108161'

# Row 161
$ws.Range("A161").Value = 'Sadriddinova Shahribonu Qoldosh qizi'
$ws.Range("B161").Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 864 soatlik'
$ws.Range("C161").Value = 'AB0664921'
$ws.Range("D161").Value = '''754'
$ws.Range("E161").Value = 'Navoiy viloyati'
$ws.Range("F161").Value = 'Nurota tumani'
$ws.Range("G161").Value = '''998941471804'
$ws.Range("H161").Value = '19-11-2024'
$ws.Range("I161").Value = '''+998942542111'

# Row 162
$ws.Range("A162").Value = 'Quvondiqova Huriyat Bahronovna'
$ws.Range("B162").Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C162").Value = 'AD2409734'
$ws.Range("D162").Value = '''755'
$ws.Range("E162").Value = 'Toshkent shahri'
$ws.Range("F162").Value = 'Mirobod tumani'
$ws.Range("G162").Value = '''+998950038686'
$ws.Range("H162").Value = '19-11-2024'
$ws.Range("I162").Value = '''+998950038686'

# Row 163
$ws.Range("A163").Value = 'TO''XSANOVA FERUZA SOBIROVNA'
$ws.Range("B163").Value = 'Defektologiya (logopediya) 576 soatlik'
$ws.Range("C163").Value = 'AD7550400'
$ws.Range("D163").Value = '''756'
$ws.Range("E163").Value = 'Navoiy viloyati'
$ws.Range("F163").Value = 'Qiziltepa tumani'
$ws.Range("G163").Value = '''998912508529'
$ws.Range("H163").Value = '19-11-2024'
$ws.Range("I163").Value = '''+998912508529'

# Row 164
$ws.Range("A164").Value = 'Saparbayeva Xurshidaxon Dostonbek qizi'
$ws.Range("B164").Value = 'Defektologiya (logopediya) 576 soatlik'
$ws.Range("C164").Value = 'AD0018405'
$ws.Range("D164").Value = '''757'
$ws.Range("E164").Value = 'Andijon viloyati'
$ws.Range("F164").Value = 'Andijon tuman'
$ws.Range("G164").Value = '''998905284303'
$ws.Range("H164").Value = '20-11-2024'
$ws.Range("I164").Value = '''+998905284303'

# Row 165
$ws.Range("A165").Value = 'Choriyeva Dilafruz Shuhrat qizi'
$ws.Range("B165").Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 864 soatlik'
$ws.Range("C165").Value = 'AD6113797'
$ws.Range("D165").Value = '''758'
$ws.Range("E165").Value = 'Qashqadaryo viloyati'
$ws.Range("F165").Value = 'Qarshi tumani'
$ws.Range("G165").Value = '''998935400332'
$ws.Range("H165").Value = '21-11-2024'
$ws.Range("I165").Value = '''+998507200332'

# Row 166
$ws.Range("A166").Value = 'Shukurlayeva Sabohat Nurlayevna'
$ws.Range("B166").Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C166").Value = 'AD7962143'
$ws.Range("D166").Value = '''759'
$ws.Range("E166").Value = 'Xorazm viloyati'
$ws.Range("F166").Value = 'Gurlan tumani'
$ws.Range("G166").Value = '''998993486603'
$ws.Range("H166").Value = '21-11-2024'
$ws.Range("I166").Value = '''+998993486603'

# Row 167
$ws.Range("A167").Value = 'Yuldasheva Irodaxon Raimqul qizi'
$ws.Range("B167").Value = 'Amaliy psixologiya 576 soatlik'
$ws.Range("C167").Value = 'AC2714983'
$ws.Range("D167").Value = '''760'
$ws.Range("E167").Value = 'Qashqadaryo viloyati'
$ws.Range("F167").Value = 'Yakkabogʻ tumani'
$ws.Range("G167").Value = '''998973095203'
$ws.Range("H167").Value = '21-11-2024'
$ws.Range("I167").Value = '''+998973095203'

# Row 168
$ws.Range("A168").Value = 'Sotvoldiyeva Nozima Qodirjonovna'
$ws.Range("B168").Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 576 soatlik'
$ws.Range("C168").Value = 'AD4032668'
$ws.Range("D168").Value = '''761'
$ws.Range("E168").Value = 'Fargona viloyati'
$ws.Range("F168").Value = 'Bagʻdod tumani'
$ws.Range("G168").Value = '''998917373670'
$ws.Range("H168").Value = '21-11-2024'
$ws.Range("I168").Value = '''+998917373670'

# Row 169
$ws.Range("A169").Value = 'Buranova Shaxnoza Olimovna'
$ws.Range("B169").Value = 'Defektologiya (logopediya) 576 soatlik'
$ws.Range("C169").Value = 'AD3858103'
$ws.Range("D169").Value = '''762'
$ws.Range("E169").Value = 'Jizzax viloyati'
$ws.Range("F169").Value = 'Sharof Rashidov tumani'
$ws.Range("G169").Value = '''998902977667'
$ws.Range("H169").Value = '21-11-2024'
$ws.Range("I169").Value = '''+998940687667'

# Row 170
$ws.Range("A170").Value = 'Safarova Madina  Baxtiyor qizi'
$ws.Range("B170").Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi 864 soatlik'
$ws.Range("C170").Value = 'AB8322263'
$ws.Range("D170").Value = '''763'
$ws.Range("E170").Value = 'Surxondaryo viloyati'
$ws.Range("F170").Value = 'Denov tumani'
$ws.Range("G170").Value = '''998975342848'
$ws.Range("H170").Value = '21-11-2024'
$ws.Range("I170").Value = '''+998975342847'

